$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating point drift on the existing A3 timestamp
$ws.Range("A3").Value = 45869.45852532407

# Append the new row of sensor data collected by the scheduled task
$ws.Range("A4").Value = 45869.50021532802
$ws.Range("B4").Value = 2025
$ws.Range("C4").Value = 31
$ws.Range("D4").Value = 22.22
$ws.Range("E4").Value = 69.77
$ws.Range("F4").Value = 632.8099999999999
$ws.Range("G4").Value = 12.15
$ws.Range("H4").Value = "ESE"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "12:00:18"

# Match the numeric date/time style used by the other rows in column A
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat
